# Change in unit of AIC and run of new results
#
# Across every yearly worksheet in the workbook, the values stored in the
# "AIC" style block (rows 5, 7 and 8 of columns D:G) were rescaled by a
# factor of 1e-6 (a unit change), together with a refreshed model run.
# Zero-valued cells are left untouched (0 * 1e-6 is still 0), matching the
# original diff exactly.

$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count

$targetCells = @("D5", "E5", "F5", "G5", "D7", "E7", "F7", "G7", "D8", "E8", "F8", "G8")
$factor = 0.000001

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)
    foreach ($addr in $targetCells) {
        $rng = $ws.Range($addr)
        $v = $rng.Value2
        if ($v -ne 0) {
            $rng.Value2 = $v * $factor
        }
    }
}
